$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Professional summary paragraph: neutralize the "affecting all
#    Black and Asian-American voters" language to "affecting 50M voters"
# ------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute(
    "Discovered systematic demographic coding errors affecting all Black and Asian-American voters, developed",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Discovered systematic demographic coding errors affecting 50M voters, developed",
    2) | Out-Null

# ------------------------------------------------------------------
# 2. Work experience bullet under "Partner - Siege Analytics": same
#    language change, but "50M" needs to become its own bold run
#    (matching the formatting already used for the "23%"/"64%" runs).
# ------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute(
    "Discovered systematic race coding errors affecting all Black and Asian-American voters, developed",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Discovered systematic race coding errors affecting 50M voters, developed",
    2) | Out-Null

$r2b = $d.Content
$r2b.Find.Execute("50M voters, developed geospatial machine learning") | Out-Null
$numRange = $d.Range($r2b.Start, $r2b.Start + 3)
$numRange.Font.Bold = 1
$numRange.Font.Color = 5258796   # RGB(0x2C,0x3E,0x50)

# ------------------------------------------------------------------
# 3. Reorder "Field Director - The Feldman Group" experience block so
#    it appears right after "Research Director - PCCC" (immediately
#    before "Software Engineer - Salsa Labs") instead of after
#    "Programmer - Lake Research Partners".
# ------------------------------------------------------------------
$feldmanStart = -1
$feldmanEnd = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt -match "Field Director - The Feldman Group") {
        $feldmanStart = $i
        $feldmanEnd = $i + 4
    }
}

$startP = $d.Paragraphs.Item($feldmanStart)
$endP = $d.Paragraphs.Item($feldmanEnd)
$moveRange = $d.Range($startP.Range.Start, $endP.Range.End)
$moveRange.Cut() | Out-Null

$salsaIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt -match "Software Engineer - Salsa Labs") {
        $salsaIdx = $i
    }
}
$prevP = $d.Paragraphs.Item($salsaIdx - 1)
$insertPoint = $d.Range($prevP.Range.End, $prevP.Range.End)
$insertPoint.Paste() | Out-Null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt -match "Field Director - The Feldman Group") {
        $d.Paragraphs.Item($i).Style = "Heading 3"
    }
}

# ------------------------------------------------------------------
# 4. "Geospatial Demographic Classification System" project Impact
#    line: same neutral-language change, with "nationwide" appended.
# ------------------------------------------------------------------
$r4 = $d.Content
$r4.Find.Execute(
    "Impact: Corrected demographic data affecting all Black and Asian-American voters, improved electoral",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Impact: Corrected demographic data affecting 50M voters nationwide, improved electoral",
    2) | Out-Null
